$d = $word.ActiveDocument

# Disable smart quotes auto-formatting so we can insert straight quotes reliably.
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false

function Split-RunAt($pos) {
    # Splits the run containing absolute character position $pos into two
    # separate <w:r> elements at that position, without altering any
    # paragraph-level identity/attributes. Achieved by temporarily
    # inserting a paragraph break at $pos (which naturally creates a new
    # run boundary), copying the tail text back onto the original
    # paragraph, then deleting the now-duplicated temporary paragraph
    # (text + its paragraph mark) so the original paragraph (and its
    # rsid/paraId attributes) survive untouched.
    $paraIdx = -1
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $pr = $d.Paragraphs.Item($i).Range
        if ($pos -ge $pr.Start -and $pos -lt $pr.End) {
            $paraIdx = $i
            break
        }
    }
    if ($paraIdx -eq -1) { throw "pos $pos not found in any paragraph" }

    $ins = $d.Range($pos, $pos)
    $ins.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($paraIdx + 1).Range
    $tailText = $d.Range($newPara.Start, $newPara.End - 1).Text

    $para2 = $d.Paragraphs.Item($paraIdx).Range
    $insertionPoint = $d.Range($para2.End - 1, $para2.End - 1)
    $insertionPoint.InsertAfter($tailText)

    $newPara2 = $d.Paragraphs.Item($paraIdx + 1).Range
    $delRange = $d.Range($newPara2.Start, $newPara2.End)
    $delRange.Delete()
}

function Move-GoBackBookmark-ToDocStart() {
    # Relocates the _GoBack bookmark from its current location (right
    # after "Bluetooth") to the very start of the document (right after
    # the Heading1 paragraph's <w:pPr>, before the "Signatures" run).
    # A direct Range(0,0) collapsed bookmark at doc-start is unreliable,
    # so we temporarily insert two placeholder characters at position 0,
    # add the bookmark right after them, then delete the placeholder.
    $old = $d.Bookmarks.Item("_GoBack")
    $old.Delete()

    $head = $d.Range(0, 0)
    $head.InsertBefore("ZZ")

    $bmRange = $d.Range(2, 2)
    $d.Bookmarks.Add("_GoBack", $bmRange)

    $placeholder = $d.Range(0, 2)
    $placeholder.Delete()
}

Move-GoBackBookmark-ToDocStart

# Locate the paragraph containing "Labs are marked as..." sentence and
# compute absolute character offsets for the curly-quote boundaries.
# (Word's Find treats straight/curly quotes as equivalent, so we scan
# character-by-character looking for the actual Unicode curly-quote
# code points instead of relying on Find/Replace here.)
$para2 = $d.Paragraphs.Item(2).Range
$openQuote = [int][char]0x201C
$closeQuote = [int][char]0x201D

$quotePositions = New-Object System.Collections.ArrayList
for ($i = $para2.Start; $i -lt $para2.End; $i++) {
    $ch = $d.Range($i, $i + 1).Text
    $code = [int][char]$ch[0]
    if ($code -eq $openQuote -or $code -eq $closeQuote) {
        [void]$quotePositions.Add($i)
    }
}

$q1start = $quotePositions[0]   # opening quote before Basic
$q2start = $quotePositions[1]   # closing quote after Basic
$q3start = $quotePositions[2]   # opening quote before Advanced
$q4start = $quotePositions[3]   # closing quote after Advanced

Write-Host "q1=$q1start q2=$q2start q3=$q3start q4=$q4start"

# Split the run into 9 pieces at the 8 boundary positions (ascending order).
$positions = @($q1start, $q1start + 1, $q2start, $q2start + 1, $q3start, $q3start + 1, $q4start, $q4start + 1)
foreach ($p in $positions) {
    Split-RunAt $p
}

# Replace each lone curly quote character with a straight quote.
$d.Range($q1start, $q1start + 1).Text = [string][char]34
$d.Range($q2start, $q2start + 1).Text = [string][char]34
$d.Range($q3start, $q3start + 1).Text = [string][char]34
$d.Range($q4start, $q4start + 1).Text = [string][char]34

Write-Host "DONE"
